# Updated cryptos list - refresh Price (column D) and Volume(1h) (column E)
# values for rows 2-51 on the active worksheet.
#
# Note: several "Price" values look like plain decimal numbers (e.g. "218.11").
# Excel's COM layer would otherwise silently coerce such inline text into a
# numeric cell value. To keep them as text (matching the source data, which
# also contains un-parseable multi-dot strings like "26.301.23" in the same
# column) we prefix the genuinely numeric-looking ones with a leading
# apostrophe, Excel's standard "force text" quote-prefix convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => @{ D = <new Price text, or $null to leave unchanged>; E = <new Volume(1h) text> }
$rowUpdates = [ordered]@{
    2  = @{ D = "26.301.23";     E = "  +0.75%  " }
    3  = @{ D = "1.679.05";      E = "  +0.70%  " }
    4  = @{ D = $null;           E = "  +0.32%  " }
    5  = @{ D = "'218.11";       E = "  +0.60%  " }
    6  = @{ D = "'0.5260";       E = "  +2.93%  " }
    7  = @{ D = "'1.008";        E = "  +0.29%  " }
    8  = @{ D = $null;           E = "  +2.18%  " }
    9  = @{ D = "'0.06464";      E = "  +0.57%  " }
    10 = @{ D = $null;           E = "  +1.13%  " }
    11 = @{ D = "'0.07510";      E = "  +1.40%  " }
    12 = @{ D = "1.677.82";      E = "  +0.62%  " }
    13 = @{ D = "'4.519";        E = "  +0.25%  " }
    14 = @{ D = "'0.5789";       E = "  -0.46%  " }
    15 = @{ D = "'0.000008504";  E = "  -0.79%  " }
    16 = @{ D = "'64.77";        E = "  +0.51%  " }
    17 = @{ D = "26.347.05";     E = "  +0.73%  " }
    18 = @{ D = "'4.924";        E = "  -0.05%  " }
    19 = @{ D = $null;           E = "  +0.28%  " }
    20 = @{ D = $null;           E = "  +0.66%  " }
    21 = @{ D = "'189.75";       E = "  +0.36%  " }
    22 = @{ D = "'6.201";        E = "  -0.21%  " }
    23 = @{ D = $null;           E = "  +0.27%  " }
    24 = @{ D = "'144.90";       E = "  -0.60%  " }
    25 = @{ D = "'7.789";        E = "  +1.98%  " }
    26 = @{ D = "'0.1257";       E = "  +5.01%  " }
    27 = @{ D = "'15.78";        E = "  +1.02%  " }
    28 = @{ D = $null;           E = "  +0.63%  " }
    29 = @{ D = $null;           E = "  +4.41%  " }
    30 = @{ D = "'1.325";        E = "  +0.45%  " }
    31 = @{ D = "'3.592";        E = "  +2.01%  " }
    32 = @{ D = "'3.588";        E = "  +2.05%  " }
    33 = @{ D = "'1.659";        E = "  +1.36%  " }
    34 = @{ D = "'1.027";        E = "  +0.97%  " }
    35 = @{ D = "'0.6215";       E = "  +2.25%  " }
    36 = @{ D = "'2.405";        E = "  +1.81%  " }
    37 = @{ D = "'2.736";        E = "  +2.52%  " }
    38 = @{ D = "'6.301";        E = "  +1.54%  " }
    39 = @{ D = "1.113.37";      E = "  +3.24%  " }
    40 = @{ D = "'0.01618";      E = "  +0.42%  " }
    41 = @{ D = "'0.8734";       E = "  +1.32%  " }
    42 = @{ D = $null;           E = "  +0.65%  " }
    43 = @{ D = "'100.53";       E = "  -0.23%  " }
    44 = @{ D = "1.829.03";      E = "  +0.74%  " }
    45 = @{ D = $null;           E = "  -2.69%  " }
    46 = @{ D = "'56.86";        E = "  +1.18%  " }
    47 = @{ D = "'8.185";        E = "  +1.51%  " }
    48 = @{ D = "'1.001";        E = "  -0.20%  " }
    49 = @{ D = "'0.05275";      E = "  +1.31%  " }
    50 = @{ D = "'0.4296";       E = "  +0.16%  " }
    51 = @{ D = $null;           E = "  +1.95%  " }
}

foreach ($row in $rowUpdates.Keys) {
    $update = $rowUpdates[$row]
    if ($null -ne $update.D) {
        $ws.Cells.Item($row, 4).Value = $update.D
    }
    $ws.Cells.Item($row, 5).Value = $update.E
}
